$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 (names shift by 2 positions because two new
#     line entries, line7/line8, are inserted before the extrN entries;
#     numeric C/D/E values are also recomputed) ---

# Row 8: line7
$ws.Cells.Item(8, 2).Value2 = "line7"
$ws.Cells.Item(8, 3).Value2 = 14
$ws.Cells.Item(8, 4).Value2 = 11
$ws.Cells.Item(8, 5).Value2 = $true

# Row 9: line8
$ws.Cells.Item(9, 2).Value2 = "line8"
$ws.Cells.Item(9, 3).Value2 = 16
$ws.Cells.Item(9, 4).Value2 = 9
$ws.Cells.Item(9, 5).Value2 = $true

# Row 10: extr1
$ws.Cells.Item(10, 2).Value2 = "extr1"
$ws.Cells.Item(10, 3).Value2 = 5
$ws.Cells.Item(10, 4).Value2 = 12
$ws.Cells.Item(10, 5).Value2 = $true

# Row 11: extr2
$ws.Cells.Item(11, 2).Value2 = "extr2"
$ws.Cells.Item(11, 3).Value2 = 5
$ws.Cells.Item(11, 4).Value2 = 9
$ws.Cells.Item(11, 5).Value2 = $true

# Row 12: extr3
$ws.Cells.Item(12, 2).Value2 = "extr3"
$ws.Cells.Item(12, 3).Value2 = 10
$ws.Cells.Item(12, 4).Value2 = 11
$ws.Cells.Item(12, 5).Value2 = $true

# Row 13: extr4
$ws.Cells.Item(13, 2).Value2 = "extr4"
$ws.Cells.Item(13, 3).Value2 = 7
$ws.Cells.Item(13, 4).Value2 = 8
$ws.Cells.Item(13, 5).Value2 = $false

# Row 14: extr5
$ws.Cells.Item(14, 2).Value2 = "extr5"
$ws.Cells.Item(14, 3).Value2 = 9
$ws.Cells.Item(14, 4).Value2 = 11
$ws.Cells.Item(14, 5).Value2 = $false

# Row 15: extr6
$ws.Cells.Item(15, 2).Value2 = "extr6"
$ws.Cells.Item(15, 3).Value2 = 7
$ws.Cells.Item(15, 4).Value2 = 11
$ws.Cells.Item(15, 5).Value2 = $true

# --- Add two new rows (16, 17) for extr7 and extr8 ---
# Copy the formatting of column A from the last existing data row so the
# new index cells keep the bold/centered/bordered style.
$ws.Cells.Item(15, 1).Copy() | Out-Null
$ws.Cells.Item(16, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15, 1).Copy() | Out-Null
$ws.Cells.Item(17, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 16: extr7
$ws.Cells.Item(16, 1).Value2 = 14
$ws.Cells.Item(16, 2).Value2 = "extr7"
$ws.Cells.Item(16, 3).Value2 = 5
$ws.Cells.Item(16, 4).Value2 = 7
$ws.Cells.Item(16, 5).Value2 = $false

# Row 17: extr8
$ws.Cells.Item(17, 1).Value2 = 15
$ws.Cells.Item(17, 2).Value2 = "extr8"
$ws.Cells.Item(17, 3).Value2 = 8
$ws.Cells.Item(17, 4).Value2 = 5
$ws.Cells.Item(17, 5).Value2 = $true
